{"js": "// Auto-generated: replace the 100 addition/subtraction equations in the\n// 20x5 table with their updated values, matching the target diff.\n// Each entry is [rowIndex, colIndex, oldText, newText] in document order.\nconst replacements = [[0, 0, \"42+32=74\", \"87-22=65\"], [0, 1, \"14+17=31\", \"79+6=85\"], [0, 2, \"57-54=3\", \"39+36=75\"], [0, 3, \"86-14=72\", \"60-8=52\"], [0, 4, \"50-46=4\", \"59-20=39\"], [1, 0, \"23+32=55\", \"92-17=75\"], [1, 1, \"31+7=38\", \"78+14=92\"], [1, 2, \"66-53=13\", \"15-3=12\"], [1, 3, \"47-9=38\", \"27+7=34\"], [1, 4, \"94-11=83\", \"58-43=15\"], [2, 0, \"87-72=15\", \"7+5=12\"], [2, 1, \"76-23=53\", \"96-2=94\"], [2, 2, \"4+19=23\", \"52-10=42\"], [2, 3, \"89-24=65\", \"68-31=37\"], [2, 4, \"13+85=98\", \"70-67=3\"], [3, 0, \"54-13=41\", \"19-10=9\"], [3, 1, \"47+8=55\", \"40+40=80\"], [3, 2, \"46+34=80\", \"44+18=62\"], [3, 3, \"64-33=31\", \"91-67=24\"], [3, 4, \"76-5=71\", \"43-42=1\"], [4, 0, \"23+74=97\", \"35-7=28\"], [4, 1, \"81-27=54\", \"29-16=13\"], [4, 2, \"81-8=73\", \"69-52=17\"], [4, 3, \"30+42=72\", \"39+26=65\"], [4, 4, \"50-11=39\", \"33-21=12\"], [5, 0, \"50-15=35\", \"76-50=26\"], [5, 1, \"38+38=76\", \"64-3=61\"], [5, 2, \"96-3=93\", \"19+2=21\"], [5, 3, \"53-1=52\", \"65+5=70\"], [5, 4, \"20+53=73\", \"22+32=54\"], [6, 0, \"81-62=19\", \"48+6=54\"], [6, 1, \"95-59=36\", \"50-13=37\"], [6, 2, \"66-3=63\", \"96-68=28\"], [6, 3, \"82-25=57\", \"11+23=34\"], [6, 4, \"81+2=83\", \"8+51=59\"], [7, 0, \"82-82=0\", \"27+70=97\"], [7, 1, \"12+57=69\", \"45+22=67\"], [7, 2, \"42-35=7\", \"64+3=67\"], [7, 3, \"82+15=97\", \"35-15=20\"], [7, 4, \"26-0=26\", \"98-82=16\"], [8, 0, \"50-49=1\", \"27+35=62\"], [8, 1, \"44+46=90\", \"23+62=85\"], [8, 2, \"24+22=46\", \"44-29=15\"], [8, 3, \"52-2=50\", \"31+64=95\"], [8, 4, \"22+2=24\", \"19+46=65\"], [9, 0, \"11+72=83\", \"55+39=94\"], [9, 1, \"82-35=47\", \"4+88=92\"], [9, 2, \"50+35=85\", \"35-27=8\"], [9, 3, \"4+75=79\", \"80-17=63\"], [9, 4, \"63-60=3\", \"50-20=30\"], [10, 0, \"10+55=65\", \"36-32=4\"], [10, 1, \"62-2=60\", \"47+23=70\"], [10, 2, \"30-1=29\", \"15+70=85\"], [10, 3, \"40+53=93\", \"20+54=74\"], [10, 4, \"53+23=76\", \"17-3=14\"], [11, 0, \"83-66=17\", \"53+6=59\"], [11, 1, \"2+43=45\", \"48-20=28\"], [11, 2, \"77+0=77\", \"61-4=57\"], [11, 3, \"51+40=91\", \"44+51=95\"], [11, 4, \"64+9=73\", \"98-0=98\"], [12, 0, \"16+82=98\", \"4+68=72\"], [12, 1, \"89-64=25\", \"2+75=77\"], [12, 2, \"53-45=8\", \"11-4=7\"], [12, 3, \"83-47=36\", \"70+17=87\"], [12, 4, \"84-30=54\", \"95-1=94\"], [13, 0, \"52-21=31\", \"61-0=61\"], [13, 1, \"57-23=34\", \"42+48=90\"], [13, 2, \"68-54=14\", \"92-3=89\"], [13, 3, \"37+31=68\", \"15-11=4\"], [13, 4, \"70-15=55\", \"95+3=98\"], [14, 0, \"66+20=86\", \"46-20=26\"], [14, 1, \"19+41=60\", \"49-22=27\"], [14, 2, \"5+68=73\", \"89-7=82\"], [14, 3, \"13+51=64\", \"79-25=54\"], [14, 4, \"78-23=55\", \"42+43=85\"], [15, 0, \"34+58=92\", \"75-30=45\"], [15, 1, \"23+12=35\", \"28+54=82\"], [15, 2, \"37-1=36\", \"48-43=5\"], [15, 3, \"22+39=61\", \"20+40=60\"], [15, 4, \"10+15=25\", \"86-67=19\"], [16, 0, \"77-1=76\", \"79-67=12\"], [16, 1, \"36+5=41\", \"86-48=38\"], [16, 2, \"81-76=5\", \"24+26=50\"], [16, 3, \"5+61=66\", \"23+71=94\"], [16, 4, \"10+87=97\", \"18+26=44\"], [17, 0, \"38-20=18\", \"21+10=31\"], [17, 1, \"19+9=28\", \"28+49=77\"], [17, 2, \"24-13=11\", \"76+19=95\"], [17, 3, \"18-12=6\", \"11+59=70\"], [17, 4, \"34-0=34\", \"19+46=65\"], [18, 0, \"57+40=97\", \"93-34=59\"], [18, 1, \"71-68=3\", \"83-38=45\"], [18, 2, \"8+48=56\", \"31-18=13\"], [18, 3, \"40+27=67\", \"82-1=81\"], [18, 4, \"80+9=89\", \"73-69=4\"], [19, 0, \"84-76=8\", \"63-47=16\"], [19, 1, \"32+8=40\", \"7+84=91\"], [19, 2, \"94-76=18\", \"30-11=19\"], [19, 3, \"77-16=61\", \"3+36=39\"], [19, 4, \"43+37=80\", \"11+23=34\"]];\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table in the document body, found none.\");\n}\nconst table = tables.items[0];\n\n// Load all cells' paragraphs in one go for efficiency.\nconst cells = [];\nfor (const [row, col] of replacements) {\n  const cell = table.getCell(row, col);\n  cell.body.paragraphs.load(\"items/text\");\n  cells.push(cell);\n}\nawait context.sync();\n\n// Sanity-check the existing text matches what the diff expects before we\n// touch anything, then replace the run's text in place (via the paragraph\n// range) so the surrounding run/paragraph formatting is preserved exactly.\nfor (let i = 0; i < replacements.length; i++) {\n  const [row, col, oldText] = replacements[i];\n  const paragraph = cells[i].body.paragraphs.items[0];\n  const actual = paragraph.text;\n  if (actual !== oldText) {\n    throw new Error(\n      `Cell (${row},${col}) expected \"${oldText}\" but found \"${actual}\"`\n    );\n  }\n}\n\nfor (let i = 0; i < replacements.length; i++) {\n  const [row, col, oldText, newText] = replacements[i];\n  const cell = cells[i];\n  const paragraph = cell.body.paragraphs.items[0];\n  const range = paragraph.getRange();\n  range.insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Auto-generated: replace the 100 addition/subtraction equations in the\n# 20x5 table with their updated values, matching the target diff.\n# Cell() is 1-indexed (row, column); Range.Text assignment replaces the\n# cell's text content while leaving the existing run/paragraph formatting\n# (font, size, alignment) untouched.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfunction Get-CellText($cell) {\n    # Cell ranges end with a cell-mark (CR + BEL); strip it for comparison.\n    return $cell.Range.Text.TrimEnd([char]13, [char]7)\n}\n\n$cell = $t.Cell(1,1)\n$current = Get-CellText $cell\nif ($current -ne \"42+32=74\") {\n    throw \"Cell (1,1) expected `\"42+32=74`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"87-22=65\"\n\n$cell = $t.Cell(1,2)\n$current = Get-CellText $cell\nif ($current -ne \"14+17=31\") {\n    throw \"Cell (1,2) expected `\"14+17=31`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"79+6=85\"\n\n$cell = $t.Cell(1,3)\n$current = Get-CellText $cell\nif ($current -ne \"57-54=3\") {\n    throw \"Cell (1,3) expected `\"57-54=3`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"39+36=75\"\n\n$cell = $t.Cell(1,4)\n$current = Get-CellText $cell\nif ($current -ne \"86-14=72\") {\n    throw \"Cell (1,4) expected `\"86-14=72`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"60-8=52\"\n\n$cell = $t.Cell(1,5)\n$current = Get-CellText $cell\nif ($current -ne \"50-46=4\") {\n    throw \"Cell (1,5) expected `\"50-46=4`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"59-20=39\"\n\n$cell = $t.Cell(2,1)\n$current = Get-CellText $cell\nif ($current -ne \"23+32=55\") {\n    throw \"Cell (2,1) expected `\"23+32=55`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"92-17=75\"\n\n$cell = $t.Cell(2,2)\n$current = Get-CellText $cell\nif ($current -ne \"31+7=38\") {\n    throw \"Cell (2,2) expected `\"31+7=38`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"78+14=92\"\n\n$cell = $t.Cell(2,3)\n$current = Get-CellText $cell\nif ($current -ne \"66-53=13\") {\n    throw \"Cell (2,3) expected `\"66-53=13`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"15-3=12\"\n\n$cell = $t.Cell(2,4)\n$current = Get-CellText $cell\nif ($current -ne \"47-9=38\") {\n    throw \"Cell (2,4) expected `\"47-9=38`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"27+7=34\"\n\n$cell = $t.Cell(2,5)\n$current = Get-CellText $cell\nif ($current -ne \"94-11=83\") {\n    throw \"Cell (2,5) expected `\"94-11=83`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"58-43=15\"\n\n$cell = $t.Cell(3,1)\n$current = Get-CellText $cell\nif ($current -ne \"87-72=15\") {\n    throw \"Cell (3,1) expected `\"87-72=15`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"7+5=12\"\n\n$cell = $t.Cell(3,2)\n$current = Get-CellText $cell\nif ($current -ne \"76-23=53\") {\n    throw \"Cell (3,2) expected `\"76-23=53`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"96-2=94\"\n\n$cell = $t.Cell(3,3)\n$current = Get-CellText $cell\nif ($current -ne \"4+19=23\") {\n    throw \"Cell (3,3) expected `\"4+19=23`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"52-10=42\"\n\n$cell = $t.Cell(3,4)\n$current = Get-CellText $cell\nif ($current -ne \"89-24=65\") {\n    throw \"Cell (3,4) expected `\"89-24=65`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"68-31=37\"\n\n$cell = $t.Cell(3,5)\n$current = Get-CellText $cell\nif ($current -ne \"13+85=98\") {\n    throw \"Cell (3,5) expected `\"13+85=98`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"70-67=3\"\n\n$cell = $t.Cell(4,1)\n$current = Get-CellText $cell\nif ($current -ne \"54-13=41\") {\n    throw \"Cell (4,1) expected `\"54-13=41`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"19-10=9\"\n\n$cell = $t.Cell(4,2)\n$current = Get-CellText $cell\nif ($current -ne \"47+8=55\") {\n    throw \"Cell (4,2) expected `\"47+8=55`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"40+40=80\"\n\n$cell = $t.Cell(4,3)\n$current = Get-CellText $cell\nif ($current -ne \"46+34=80\") {\n    throw \"Cell (4,3) expected `\"46+34=80`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"44+18=62\"\n\n$cell = $t.Cell(4,4)\n$current = Get-CellText $cell\nif ($current -ne \"64-33=31\") {\n    throw \"Cell (4,4) expected `\"64-33=31`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"91-67=24\"\n\n$cell = $t.Cell(4,5)\n$current = Get-CellText $cell\nif ($current -ne \"76-5=71\") {\n    throw \"Cell (4,5) expected `\"76-5=71`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"43-42=1\"\n\n$cell = $t.Cell(5,1)\n$current = Get-CellText $cell\nif ($current -ne \"23+74=97\") {\n    throw \"Cell (5,1) expected `\"23+74=97`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"35-7=28\"\n\n$cell = $t.Cell(5,2)\n$current = Get-CellText $cell\nif ($current -ne \"81-27=54\") {\n    throw \"Cell (5,2) expected `\"81-27=54`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"29-16=13\"\n\n$cell = $t.Cell(5,3)\n$current = Get-CellText $cell\nif ($current -ne \"81-8=73\") {\n    throw \"Cell (5,3) expected `\"81-8=73`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"69-52=17\"\n\n$cell = $t.Cell(5,4)\n$current = Get-CellText $cell\nif ($current -ne \"30+42=72\") {\n    throw \"Cell (5,4) expected `\"30+42=72`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"39+26=65\"\n\n$cell = $t.Cell(5,5)\n$current = Get-CellText $cell\nif ($current -ne \"50-11=39\") {\n    throw \"Cell (5,5) expected `\"50-11=39`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"33-21=12\"\n\n$cell = $t.Cell(6,1)\n$current = Get-CellText $cell\nif ($current -ne \"50-15=35\") {\n    throw \"Cell (6,1) expected `\"50-15=35`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"76-50=26\"\n\n$cell = $t.Cell(6,2)\n$current = Get-CellText $cell\nif ($current -ne \"38+38=76\") {\n    throw \"Cell (6,2) expected `\"38+38=76`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"64-3=61\"\n\n$cell = $t.Cell(6,3)\n$current = Get-CellText $cell\nif ($current -ne \"96-3=93\") {\n    throw \"Cell (6,3) expected `\"96-3=93`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"19+2=21\"\n\n$cell = $t.Cell(6,4)\n$current = Get-CellText $cell\nif ($current -ne \"53-1=52\") {\n    throw \"Cell (6,4) expected `\"53-1=52`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"65+5=70\"\n\n$cell = $t.Cell(6,5)\n$current = Get-CellText $cell\nif ($current -ne \"20+53=73\") {\n    throw \"Cell (6,5) expected `\"20+53=73`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"22+32=54\"\n\n$cell = $t.Cell(7,1)\n$current = Get-CellText $cell\nif ($current -ne \"81-62=19\") {\n    throw \"Cell (7,1) expected `\"81-62=19`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"48+6=54\"\n\n$cell = $t.Cell(7,2)\n$current = Get-CellText $cell\nif ($current -ne \"95-59=36\") {\n    throw \"Cell (7,2) expected `\"95-59=36`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"50-13=37\"\n\n$cell = $t.Cell(7,3)\n$current = Get-CellText $cell\nif ($current -ne \"66-3=63\") {\n    throw \"Cell (7,3) expected `\"66-3=63`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"96-68=28\"\n\n$cell = $t.Cell(7,4)\n$current = Get-CellText $cell\nif ($current -ne \"82-25=57\") {\n    throw \"Cell (7,4) expected `\"82-25=57`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"11+23=34\"\n\n$cell = $t.Cell(7,5)\n$current = Get-CellText $cell\nif ($current -ne \"81+2=83\") {\n    throw \"Cell (7,5) expected `\"81+2=83`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"8+51=59\"\n\n$cell = $t.Cell(8,1)\n$current = Get-CellText $cell\nif ($current -ne \"82-82=0\") {\n    throw \"Cell (8,1) expected `\"82-82=0`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"27+70=97\"\n\n$cell = $t.Cell(8,2)\n$current = Get-CellText $cell\nif ($current -ne \"12+57=69\") {\n    throw \"Cell (8,2) expected `\"12+57=69`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"45+22=67\"\n\n$cell = $t.Cell(8,3)\n$current = Get-CellText $cell\nif ($current -ne \"42-35=7\") {\n    throw \"Cell (8,3) expected `\"42-35=7`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"64+3=67\"\n\n$cell = $t.Cell(8,4)\n$current = Get-CellText $cell\nif ($current -ne \"82+15=97\") {\n    throw \"Cell (8,4) expected `\"82+15=97`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"35-15=20\"\n\n$cell = $t.Cell(8,5)\n$current = Get-CellText $cell\nif ($current -ne \"26-0=26\") {\n    throw \"Cell (8,5) expected `\"26-0=26`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"98-82=16\"\n\n$cell = $t.Cell(9,1)\n$current = Get-CellText $cell\nif ($current -ne \"50-49=1\") {\n    throw \"Cell (9,1) expected `\"50-49=1`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"27+35=62\"\n\n$cell = $t.Cell(9,2)\n$current = Get-CellText $cell\nif ($current -ne \"44+46=90\") {\n    throw \"Cell (9,2) expected `\"44+46=90`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"23+62=85\"\n\n$cell = $t.Cell(9,3)\n$current = Get-CellText $cell\nif ($current -ne \"24+22=46\") {\n    throw \"Cell (9,3) expected `\"24+22=46`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"44-29=15\"\n\n$cell = $t.Cell(9,4)\n$current = Get-CellText $cell\nif ($current -ne \"52-2=50\") {\n    throw \"Cell (9,4) expected `\"52-2=50`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"31+64=95\"\n\n$cell = $t.Cell(9,5)\n$current = Get-CellText $cell\nif ($current -ne \"22+2=24\") {\n    throw \"Cell (9,5) expected `\"22+2=24`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"19+46=65\"\n\n$cell = $t.Cell(10,1)\n$current = Get-CellText $cell\nif ($current -ne \"11+72=83\") {\n    throw \"Cell (10,1) expected `\"11+72=83`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"55+39=94\"\n\n$cell = $t.Cell(10,2)\n$current = Get-CellText $cell\nif ($current -ne \"82-35=47\") {\n    throw \"Cell (10,2) expected `\"82-35=47`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"4+88=92\"\n\n$cell = $t.Cell(10,3)\n$current = Get-CellText $cell\nif ($current -ne \"50+35=85\") {\n    throw \"Cell (10,3) expected `\"50+35=85`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"35-27=8\"\n\n$cell = $t.Cell(10,4)\n$current = Get-CellText $cell\nif ($current -ne \"4+75=79\") {\n    throw \"Cell (10,4) expected `\"4+75=79`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"80-17=63\"\n\n$cell = $t.Cell(10,5)\n$current = Get-CellText $cell\nif ($current -ne \"63-60=3\") {\n    throw \"Cell (10,5) expected `\"63-60=3`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"50-20=30\"\n\n$cell = $t.Cell(11,1)\n$current = Get-CellText $cell\nif ($current -ne \"10+55=65\") {\n    throw \"Cell (11,1) expected `\"10+55=65`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"36-32=4\"\n\n$cell = $t.Cell(11,2)\n$current = Get-CellText $cell\nif ($current -ne \"62-2=60\") {\n    throw \"Cell (11,2) expected `\"62-2=60`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"47+23=70\"\n\n$cell = $t.Cell(11,3)\n$current = Get-CellText $cell\nif ($current -ne \"30-1=29\") {\n    throw \"Cell (11,3) expected `\"30-1=29`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"15+70=85\"\n\n$cell = $t.Cell(11,4)\n$current = Get-CellText $cell\nif ($current -ne \"40+53=93\") {\n    throw \"Cell (11,4) expected `\"40+53=93`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"20+54=74\"\n\n$cell = $t.Cell(11,5)\n$current = Get-CellText $cell\nif ($current -ne \"53+23=76\") {\n    throw \"Cell (11,5) expected `\"53+23=76`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"17-3=14\"\n\n$cell = $t.Cell(12,1)\n$current = Get-CellText $cell\nif ($current -ne \"83-66=17\") {\n    throw \"Cell (12,1) expected `\"83-66=17`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"53+6=59\"\n\n$cell = $t.Cell(12,2)\n$current = Get-CellText $cell\nif ($current -ne \"2+43=45\") {\n    throw \"Cell (12,2) expected `\"2+43=45`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"48-20=28\"\n\n$cell = $t.Cell(12,3)\n$current = Get-CellText $cell\nif ($current -ne \"77+0=77\") {\n    throw \"Cell (12,3) expected `\"77+0=77`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"61-4=57\"\n\n$cell = $t.Cell(12,4)\n$current = Get-CellText $cell\nif ($current -ne \"51+40=91\") {\n    throw \"Cell (12,4) expected `\"51+40=91`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"44+51=95\"\n\n$cell = $t.Cell(12,5)\n$current = Get-CellText $cell\nif ($current -ne \"64+9=73\") {\n    throw \"Cell (12,5) expected `\"64+9=73`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"98-0=98\"\n\n$cell = $t.Cell(13,1)\n$current = Get-CellText $cell\nif ($current -ne \"16+82=98\") {\n    throw \"Cell (13,1) expected `\"16+82=98`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"4+68=72\"\n\n$cell = $t.Cell(13,2)\n$current = Get-CellText $cell\nif ($current -ne \"89-64=25\") {\n    throw \"Cell (13,2) expected `\"89-64=25`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"2+75=77\"\n\n$cell = $t.Cell(13,3)\n$current = Get-CellText $cell\nif ($current -ne \"53-45=8\") {\n    throw \"Cell (13,3) expected `\"53-45=8`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"11-4=7\"\n\n$cell = $t.Cell(13,4)\n$current = Get-CellText $cell\nif ($current -ne \"83-47=36\") {\n    throw \"Cell (13,4) expected `\"83-47=36`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"70+17=87\"\n\n$cell = $t.Cell(13,5)\n$current = Get-CellText $cell\nif ($current -ne \"84-30=54\") {\n    throw \"Cell (13,5) expected `\"84-30=54`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"95-1=94\"\n\n$cell = $t.Cell(14,1)\n$current = Get-CellText $cell\nif ($current -ne \"52-21=31\") {\n    throw \"Cell (14,1) expected `\"52-21=31`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"61-0=61\"\n\n$cell = $t.Cell(14,2)\n$current = Get-CellText $cell\nif ($current -ne \"57-23=34\") {\n    throw \"Cell (14,2) expected `\"57-23=34`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"42+48=90\"\n\n$cell = $t.Cell(14,3)\n$current = Get-CellText $cell\nif ($current -ne \"68-54=14\") {\n    throw \"Cell (14,3) expected `\"68-54=14`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"92-3=89\"\n\n$cell = $t.Cell(14,4)\n$current = Get-CellText $cell\nif ($current -ne \"37+31=68\") {\n    throw \"Cell (14,4) expected `\"37+31=68`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"15-11=4\"\n\n$cell = $t.Cell(14,5)\n$current = Get-CellText $cell\nif ($current -ne \"70-15=55\") {\n    throw \"Cell (14,5) expected `\"70-15=55`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"95+3=98\"\n\n$cell = $t.Cell(15,1)\n$current = Get-CellText $cell\nif ($current -ne \"66+20=86\") {\n    throw \"Cell (15,1) expected `\"66+20=86`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"46-20=26\"\n\n$cell = $t.Cell(15,2)\n$current = Get-CellText $cell\nif ($current -ne \"19+41=60\") {\n    throw \"Cell (15,2) expected `\"19+41=60`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"49-22=27\"\n\n$cell = $t.Cell(15,3)\n$current = Get-CellText $cell\nif ($current -ne \"5+68=73\") {\n    throw \"Cell (15,3) expected `\"5+68=73`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"89-7=82\"\n\n$cell = $t.Cell(15,4)\n$current = Get-CellText $cell\nif ($current -ne \"13+51=64\") {\n    throw \"Cell (15,4) expected `\"13+51=64`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"79-25=54\"\n\n$cell = $t.Cell(15,5)\n$current = Get-CellText $cell\nif ($current -ne \"78-23=55\") {\n    throw \"Cell (15,5) expected `\"78-23=55`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"42+43=85\"\n\n$cell = $t.Cell(16,1)\n$current = Get-CellText $cell\nif ($current -ne \"34+58=92\") {\n    throw \"Cell (16,1) expected `\"34+58=92`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"75-30=45\"\n\n$cell = $t.Cell(16,2)\n$current = Get-CellText $cell\nif ($current -ne \"23+12=35\") {\n    throw \"Cell (16,2) expected `\"23+12=35`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"28+54=82\"\n\n$cell = $t.Cell(16,3)\n$current = Get-CellText $cell\nif ($current -ne \"37-1=36\") {\n    throw \"Cell (16,3) expected `\"37-1=36`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"48-43=5\"\n\n$cell = $t.Cell(16,4)\n$current = Get-CellText $cell\nif ($current -ne \"22+39=61\") {\n    throw \"Cell (16,4) expected `\"22+39=61`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"20+40=60\"\n\n$cell = $t.Cell(16,5)\n$current = Get-CellText $cell\nif ($current -ne \"10+15=25\") {\n    throw \"Cell (16,5) expected `\"10+15=25`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"86-67=19\"\n\n$cell = $t.Cell(17,1)\n$current = Get-CellText $cell\nif ($current -ne \"77-1=76\") {\n    throw \"Cell (17,1) expected `\"77-1=76`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"79-67=12\"\n\n$cell = $t.Cell(17,2)\n$current = Get-CellText $cell\nif ($current -ne \"36+5=41\") {\n    throw \"Cell (17,2) expected `\"36+5=41`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"86-48=38\"\n\n$cell = $t.Cell(17,3)\n$current = Get-CellText $cell\nif ($current -ne \"81-76=5\") {\n    throw \"Cell (17,3) expected `\"81-76=5`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"24+26=50\"\n\n$cell = $t.Cell(17,4)\n$current = Get-CellText $cell\nif ($current -ne \"5+61=66\") {\n    throw \"Cell (17,4) expected `\"5+61=66`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"23+71=94\"\n\n$cell = $t.Cell(17,5)\n$current = Get-CellText $cell\nif ($current -ne \"10+87=97\") {\n    throw \"Cell (17,5) expected `\"10+87=97`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"18+26=44\"\n\n$cell = $t.Cell(18,1)\n$current = Get-CellText $cell\nif ($current -ne \"38-20=18\") {\n    throw \"Cell (18,1) expected `\"38-20=18`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"21+10=31\"\n\n$cell = $t.Cell(18,2)\n$current = Get-CellText $cell\nif ($current -ne \"19+9=28\") {\n    throw \"Cell (18,2) expected `\"19+9=28`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"28+49=77\"\n\n$cell = $t.Cell(18,3)\n$current = Get-CellText $cell\nif ($current -ne \"24-13=11\") {\n    throw \"Cell (18,3) expected `\"24-13=11`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"76+19=95\"\n\n$cell = $t.Cell(18,4)\n$current = Get-CellText $cell\nif ($current -ne \"18-12=6\") {\n    throw \"Cell (18,4) expected `\"18-12=6`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"11+59=70\"\n\n$cell = $t.Cell(18,5)\n$current = Get-CellText $cell\nif ($current -ne \"34-0=34\") {\n    throw \"Cell (18,5) expected `\"34-0=34`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"19+46=65\"\n\n$cell = $t.Cell(19,1)\n$current = Get-CellText $cell\nif ($current -ne \"57+40=97\") {\n    throw \"Cell (19,1) expected `\"57+40=97`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"93-34=59\"\n\n$cell = $t.Cell(19,2)\n$current = Get-CellText $cell\nif ($current -ne \"71-68=3\") {\n    throw \"Cell (19,2) expected `\"71-68=3`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"83-38=45\"\n\n$cell = $t.Cell(19,3)\n$current = Get-CellText $cell\nif ($current -ne \"8+48=56\") {\n    throw \"Cell (19,3) expected `\"8+48=56`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"31-18=13\"\n\n$cell = $t.Cell(19,4)\n$current = Get-CellText $cell\nif ($current -ne \"40+27=67\") {\n    throw \"Cell (19,4) expected `\"40+27=67`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"82-1=81\"\n\n$cell = $t.Cell(19,5)\n$current = Get-CellText $cell\nif ($current -ne \"80+9=89\") {\n    throw \"Cell (19,5) expected `\"80+9=89`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"73-69=4\"\n\n$cell = $t.Cell(20,1)\n$current = Get-CellText $cell\nif ($current -ne \"84-76=8\") {\n    throw \"Cell (20,1) expected `\"84-76=8`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"63-47=16\"\n\n$cell = $t.Cell(20,2)\n$current = Get-CellText $cell\nif ($current -ne \"32+8=40\") {\n    throw \"Cell (20,2) expected `\"32+8=40`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"7+84=91\"\n\n$cell = $t.Cell(20,3)\n$current = Get-CellText $cell\nif ($current -ne \"94-76=18\") {\n    throw \"Cell (20,3) expected `\"94-76=18`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"30-11=19\"\n\n$cell = $t.Cell(20,4)\n$current = Get-CellText $cell\nif ($current -ne \"77-16=61\") {\n    throw \"Cell (20,4) expected `\"77-16=61`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"3+36=39\"\n\n$cell = $t.Cell(20,5)\n$current = Get-CellText $cell\nif ($current -ne \"43+37=80\") {\n    throw \"Cell (20,5) expected `\"43+37=80`\" but found `\"$current`\"\"\n}\n$cell.Range.Text = \"11+23=34\"\n\n"}
